$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Mark the DOCS Suite row's Sanity Runmode as "Y" (was "N") so both the
# Transmittal and DOCS suites are flagged for execution.
$ws.Range("F3").Value = "Y"
